$wb = $excel.ActiveWorkbook

# The sheet previously named "MELHORIA" is the visible/active demonstrativo
# financeiro sheet; rename it to "Demonstrativo" (the hidden "Atual" sheet
# is left untouched).
$ws = $wb.Worksheets.Item("MELHORIA")
$ws.Name = "Demonstrativo"

# Keep it the active/selected sheet (matches tabSelected="true" on this
# sheet in the original workbook).
$ws.Activate()
